# Slide 6 ("TESTS/CI/...") layout rework:
#   - remove the "Прямоугольник 6" shape (the duplicate/draft header text box)
#   - shift the remaining shapes up to close the resulting gap

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# 1) Delete the "Прямоугольник 6" shape (3rd shape in the slide's z-order).
$s.Shapes.Item(3).Delete()

# After the delete, every following shape's index shifts down by one:
#   3 -> "TextBox 1" (the tests bullet list)
#   4 -> "Рисунок 3"  (the screenshot picture)
#   5 -> "TextBox 8"  (the integration bullet list)

# 2) "TextBox 1": keep its horizontal position, move it up (y: 2079321 -> 1418067 EMU).
$textBox1 = $s.Shapes.Item(3)
$textBox1.Top = 111.65882110595703

# 3) "Рисунок 3" picture: keep its horizontal position, move it up (y: 3582370 -> 3236231 EMU).
$picture = $s.Shapes.Item(4)
$picture.Top = 254.82135009765625

# 4) "TextBox 8": move it both horizontally and vertically
#    (x: 500283 -> 569932 EMU, y: 4782878 -> 4516603 EMU).
$textBox8 = $s.Shapes.Item(5)
$textBox8.Left = 44.87653732299805
$textBox8.Top = 355.6380615234375
